# Update countries & provincias Spain
#
# The source publishes a ranking table (by total cases) that gets re-sorted
# as case counts change. Between the two snapshots a few neighbouring
# countries swapped rank/row position, and a handful of rows picked up
# refreshed totals. This script reproduces both: the row swaps (by writing
# the country name together with its refreshed stats into the right rows)
# and the plain numeric refreshes, then finally bumps the "updated at"
# timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    if ($country -ne $null) {
        $ws.Cells.Item($row, 1).Value = $country
    }
    $ws.Cells.Item($row, 2).Value = $casosTotales
    $ws.Cells.Item($row, 3).Value = $nuevosCasos
    $ws.Cells.Item($row, 4).Value = $casosActivos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $casosCriticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Row 7 - Rusia (rank 11): refreshed stats, no reorder
Set-Row 7 $null 1009995 4995 826935 165532 0 114 17528

# Rows 59-60: Kirguistan/Armenia swap rank - Armenia moves ahead of Kirguistan.
# Row 59 now carries Armenia's refreshed stats, row 60 carries Kirguistan's
# previous (row-59) stats.
Set-Row 59 "Armenia" 44271 196 38855 4529 0 3 887
Set-Row 60 "Kirguistan" 44135 99 39174 3902 0 0 1059

# Row 65 - Afganistan (rank 69): refreshed stats, no reorder
Set-Row 65 $null 38288 45 29390 7489 0 0 1409

# Rows 105-107: Luxemburgo/Zimbabue/Hungria reorder - Hungria moves ahead of
# Luxemburgo (and Zimbabue shifts down one row). Row 105 gets Hungria's
# refreshed stats; rows 106-107 inherit the previous rows' (104/105/106)
# stats, shifting the chain down by one.
Set-Row 105 "Hungria" 6923 301 3930 2373 0 1 620
Set-Row 106 "Luxemburgo" 6745 0 0 0 0 0 124
Set-Row 107 "Zimbabue" 6638 0 5250 1182 0 0 206

# Row 152 - Georgia (rank 156): refreshed stats, no reorder
Set-Row 152 $null 1568 20 1279 270 0 0 19

# Row 155 - Letonia (rank 159): refreshed stats, no reorder
Set-Row 155 $null 1410 4 1187 189 0 0 34

# Bump the "updated at" timestamp shown in row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 3 de Septiembre de 2020 a las 09:36"
